$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from F1 so G1 matches the bold/border/centered style
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Populate column G (header + 92 data rows) per the diff
$ws.Range("G1").Value = "G"
$ws.Range("G2").Value = "Hb 47"
$ws.Range("G3").Value = "Hb 48"
$ws.Range("G4").Value = "S 6"
$ws.Range("G5").Value = "Hb 7"
$ws.Range("G6").Value = "Hb 46"
$ws.Range("G7").Value = "Hb 1"
$ws.Range("G8").Value = "Hb 2"
$ws.Range("G9").Value = "Hb 3"
$ws.Range("G10").Value = "Hb 5"
$ws.Range("G11").Value = "S 24"
$ws.Range("G12").Value = "S 25"
$ws.Range("G13").Value = "S 26"
$ws.Range("G14").Value = "S 27"
$ws.Range("G15").Value = "S 28"
$ws.Range("G16").Value = "Hb 19"
$ws.Range("G17").Value = "Hb 20"
$ws.Range("G18").Value = "32 FO1Hepi"
$ws.Range("G19").Value = "33 FO2H"
$ws.Range("G20").Value = "S 29"
$ws.Range("G21").Value = "S 30"
$ws.Range("G22").Value = "S 17"
$ws.Range("G23").Value = "20 KR1Ph"
$ws.Range("G24").Value = "21 KR2T"
$ws.Range("G25").Value = "22 KR3S"
$ws.Range("G26").Value = "23 KR4V"
$ws.Range("G27").Value = "24 KR5Mt"
$ws.Range("G28").Value = "71 KR1BuTy"
$ws.Range("G29").Value = "72 KR2Sc"
$ws.Range("G30").Value = "73 KR3Mt"
$ws.Range("G31").Value = "74 KR4Fi"
$ws.Range("G32").Value = "75 KR5Mt"
$ws.Range("G33").Value = "Hb 83"
$ws.Range("G34").Value = "Hb 84"
$ws.Range("G35").Value = "Hb 85"
$ws.Range("G36").Value = "Hb 86"
$ws.Range("G37").Value = "Hb 87"
$ws.Range("G38").Value = "Hb 88"
$ws.Range("G39").Value = "Hb 89"
$ws.Range("G40").Value = "Hb 90"
$ws.Range("G41").Value = "Hb 91"
$ws.Range("G42").Value = "Hb 92"
$ws.Range("G43").Value = "Hb 50"
$ws.Range("G44").Value = "Hb 40"
$ws.Range("G45").Value = "Hb 41"
$ws.Range("G46").Value = "Hb 42"
$ws.Range("G47").Value = "Hb 43"
$ws.Range("G48").Value = "S 8"
$ws.Range("G49").Value = "S 9"
$ws.Range("G50").Value = "S 10"
$ws.Range("G51").Value = "S 11"
$ws.Range("G52").Value = "S 12"
$ws.Range("G53").Value = "Hb 53"
$ws.Range("G54").Value = "Hb 54"
$ws.Range("G55").Value = "Hb 55"
$ws.Range("G56").Value = "Hb 56"
$ws.Range("G57").Value = "Hb 57"
$ws.Range("G58").Value = "Hb 58"
$ws.Range("G59").Value = "Hb 59"
$ws.Range("G60").Value = "Hb 60"
$ws.Range("G61").Value = "Hb 61"
$ws.Range("G62").Value = "Hb 62"
$ws.Range("G63").Value = "Hb 28"
$ws.Range("G64").Value = "Hb 35"
$ws.Range("G65").Value = "Hb 36"
$ws.Range("G66").Value = "Hb 38"
$ws.Range("G67").Value = "Hb 39"
$ws.Range("G68").Value = "S 1"
$ws.Range("G69").Value = "S 2"
$ws.Range("G70").Value = "S 3"
$ws.Range("G71").Value = "S 4"
$ws.Range("G72").Value = "S 5"
$ws.Range("G73").Value = "Hb 30"
$ws.Range("G74").Value = "KHb 25"
$ws.Range("G75").Value = "26 SR2Fi"
$ws.Range("G76").Value = "Hb 73"
$ws.Range("G77").Value = "Hb 74"
$ws.Range("G78").Value = "Hb 75"
$ws.Range("G79").Value = "Hb 76"
$ws.Range("G80").Value = "Hb 77"
$ws.Range("G81").Value = "Hb 78"
$ws.Range("G82").Value = "Hb 79"
$ws.Range("G83").Value = "Hb 23"
$ws.Range("G84").Value = "Hb 24"
$ws.Range("G85").Value = "Hb 25"
$ws.Range("G86").Value = "Hb 26"
$ws.Range("G87").Value = "Hb 27"
$ws.Range("G88").Value = "KHb 30"
$ws.Range("G89").Value = "KHb 31"
$ws.Range("G90").Value = "27 SH1Ph"
$ws.Range("G91").Value = "28 SH2R"
$ws.Range("G92").Value = "29 SH3V"
$ws.Range("G93").Value = "KS 76"

$excel.CutCopyMode = 0
